$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Num($cellRange, $strVal) {
    $cellRange.Value = [double]$strVal
}

# Ensure row 5 exists with the same row-header style as rows 2-4 (copy format from row 4)
$ws.Range("A4:T4").Copy($ws.Range("A5:T5"))

# --- String-valued cells must be written in column-major order (C, then O..T), top-to-bottom, ---
# --- so new shared-string entries are appended in the same order as the target workbook.      ---
$ws.Range("C3").Value = "['N1ratio-ArgsPreds', 'latitude', 'longitude', 'Macro_class']"
$ws.Range("C4").Value = "['N1ratio-ArgsPreds', 'latitude', 'longitude', 'Macro_class', 'Fam_class']"
$ws.Range("C5").Value = "['N1ratio-ArgsPreds', 'latitude', 'longitude', 'Macro_class', 'Fam_class', 'Nlen_freq', 'Vlen_freq']"

$ws.Range("O2").Value = "{'const': 0.8158622228498638, 'N1ratio-ArgsPreds': -0.23974231781280247}"
$ws.Range("O3").Value = "{'const': 0.5707301483093369, 'N1ratio-ArgsPreds': -0.23736630054157704, 'latitude': 0.004125113705238998, 'longitude': 0.0003535003417894755, 'Macro_class': 0.07277894332771967}"
$ws.Range("O4").Value = "{'const': 0.6961464576636499, 'N1ratio-ArgsPreds': -0.234080535081691, 'latitude': 0.003932715362527422, 'longitude': -6.492392921279064e-05, 'Macro_class': 0.058273569555529196, 'Fam_class': -0.0017571297445560055}"
$ws.Range("O5").Value = "{'const': 0.3013014345823354, 'N1ratio-ArgsPreds': -0.2187528504296818, 'latitude': 0.0038626991338975582, 'longitude': -0.0005816308890785551, 'Macro_class': 0.03965047136416332, 'Fam_class': -0.004469086952158804, 'Nlen_freq': 0.047799212071237, 'Vlen_freq': 0.03496604140457689}"

$ws.Range("P2").Value = "{'const': 3.861723678139761e-68, 'N1ratio-ArgsPreds': 1.4607799151140621e-43}"
$ws.Range("P3").Value = "{'const': 4.590399335908072e-33, 'N1ratio-ArgsPreds': 6.126742452179407e-49, 'latitude': 5.991196955199697e-07, 'longitude': 0.16164137009314025, 'Macro_class': 9.880785720371535e-12}"
$ws.Range("P4").Value = "{'const': 1.2854831552068233e-13, 'N1ratio-ArgsPreds': 3.5622512818011976e-47, 'latitude': 2.3833674210851503e-06, 'longitude': 0.8596779117771504, 'Macro_class': 3.455875400858192e-05, 'Fam_class': 0.11769560292191579}"
$ws.Range("P5").Value = "{'const': 0.01167273354297994, 'N1ratio-ArgsPreds': 9.545059203789041e-42, 'latitude': 2.1179641117540256e-06, 'longitude': 0.11720723741364805, 'Macro_class': 0.005383749551475992, 'Fam_class': 0.0002506822457137312, 'Nlen_freq': 0.07177602002568405, 'Vlen_freq': 0.06612207205839887}"

$ws.Range("Q2").Value = "{'N1ratio-ArgsPreds': -0.5479438411342928}"
$ws.Range("Q3").Value = "{'N1ratio-ArgsPreds': -0.5425133270637095, 'latitude': 0.18532462755220405, 'longitude': 0.05907319878123149, 'Macro_class': 0.2765293156965464}"
$ws.Range("Q4").Value = "{'N1ratio-ArgsPreds': -0.5350035350354112, 'latitude': 0.17668095037079692, 'longitude': -0.010849393119766336, 'Macro_class': 0.22141500790721202, 'Fam_class': -0.11935728648818478}"
$ws.Range("Q5").Value = "{'N1ratio-ArgsPreds': -0.4999712950848697, 'latitude': 0.17353540520025776, 'longitude': -0.0971959375029516, 'Macro_class': 0.15065508252854065, 'Fam_class': -0.30357353709484436, 'Nlen_freq': 0.10205990338527482, 'Vlen_freq': 0.10161364849679277}"

$ws.Range("R2").Value = "{'N1ratio-ArgsPreds': -0.5479438411342936}"
$ws.Range("R3").Value = "{'N1ratio-ArgsPreds': -0.5772079889797707, 'latitude': 0.21360542253202497, 'longitude': 0.06053820186776226, 'Macro_class': 0.2884699120894178}"
$ws.Range("R4").Value = "{'N1ratio-ArgsPreds': -0.5687521701925946, 'latitude': 0.20233083487677564, 'longitude': -0.007660853127783767, 'Macro_class': 0.1780213224686194, 'Fam_class': -0.06772080829358151}"
$ws.Range("R5").Value = "{'N1ratio-ArgsPreds': -0.5404188496272267, 'latitude': 0.20371213287087744, 'longitude': -0.06793824368417971, 'Macro_class': 0.12039455798165379, 'Fam_class': -0.15797378583359561, 'Nlen_freq': 0.07805411206423138, 'Vlen_freq': 0.07965558414275081}"

$ws.Range("S2").Value = "{'N1ratio-ArgsPreds': -0.5479438411342935}"
$ws.Range("S3").Value = "{'N1ratio-ArgsPreds': -0.534463625943419, 'latitude': 0.16532820660595793, 'longitude': 0.04585856776705677, 'Macro_class': 0.2278035235834368}"
$ws.Range("S4").Value = "{'N1ratio-ArgsPreds': -0.5216485146302986, 'latitude': 0.15585983719914648, 'longitude': -0.005779435483898172, 'Macro_class': 0.13647739071120824, 'Fam_class': -0.05120540714737771}"
$ws.Range("S5").Value = "{'N1ratio-ArgsPreds': -0.4717623989193407, 'latitude': 0.15283185281610756, 'longitude': -0.05001637956245846, 'Macro_class': 0.0890780616243271, 'Fam_class': -0.11750766774176961, 'Nlen_freq': 0.057506394876072425, 'Vlen_freq': 0.05869373797906177}"

$ws.Range("T2").Value = "{'N1ratio-ArgsPreds': 30.024245303700386}"
$ws.Range("T3").Value = "{'N1ratio-ArgsPreds': 28.56513674565869, 'latitude': 2.733341589954231, 'longitude': 0.21030082376457382, 'Macro_class': 5.189444535702945}"
$ws.Range("T4").Value = "{'N1ratio-ArgsPreds': 27.211717281599686, 'latitude': 2.4292288851744446, 'longitude': 0.0033401874512541297, 'Macro_class': 1.8626078175339789, 'Fam_class': 0.26219937211287203}"
$ws.Range("T5").Value = "{'N1ratio-ArgsPreds': 22.255976103413115, 'latitude': 2.3357575235204364, 'longitude': 0.25016382245359126, 'Macro_class': 0.7934901062747418, 'Fam_class': 1.3808051978110123, 'Nlen_freq': 0.3306985451642769, 'Vlen_freq': 0.34449548779547584}"

# --- Numeric cells ---
Set-Num $ws.Range("A2") "0"
Set-Num $ws.Range("B2") "1"
Set-Num $ws.Range("D2") "539"
Set-Num $ws.Range("E2") "537"
Set-Num $ws.Range("F2") "1"
Set-Num $ws.Range("G2") "0.3002424530370029"
Set-Num $ws.Range("H2") "230.4086579424865"
Set-Num $ws.Range("I2") "1.460779915114385E-43"
Set-Num $ws.Range("J2") "71.15456703886447"
Set-Num $ws.Range("K2") "101.6846011131726"
Set-Num $ws.Range("L2") "30.53003407430809"
Set-Num $ws.Range("M2") "0.1325038492343845"
Set-Num $ws.Range("N2") "0.1890048347828486"

Set-Num $ws.Range("A3") "1"
Set-Num $ws.Range("B3") "2"
Set-Num $ws.Range("D3") "539"
Set-Num $ws.Range("E3") "534"
Set-Num $ws.Range("F3") "4"
Set-Num $ws.Range("G3") "0.4282747392666378"
Set-Num $ws.Range("H3") "100.0037633788599"
Set-Num $ws.Range("I3") "1.702338642847579E-63"
Set-Num $ws.Range("J3") "58.13565508399651"
Set-Num $ws.Range("K3") "101.6846011131726"
Set-Num $ws.Range("L3") "10.88723650729401"
Set-Num $ws.Range("M3") "0.108868267947559"
Set-Num $ws.Range("N3") "0.1890048347828486"
Set-Num $ws.Range("U3") "0.1280322862296349"
Set-Num $ws.Range("V3") "39.86136089149219"
Set-Num $ws.Range("W3") "2.945126212213881E-23"

Set-Num $ws.Range("A4") "2"
Set-Num $ws.Range("B4") "3"
Set-Num $ws.Range("D4") "539"
Set-Num $ws.Range("E4") "533"
Set-Num $ws.Range("F4") "5"
Set-Num $ws.Range("G4") "0.4308967329877673"
Set-Num $ws.Range("H4") "80.71222640777542"
Set-Num $ws.Range("I4") "5.384584677552259E-63"
Set-Num $ws.Range("J4") "57.86903869834222"
Set-Num $ws.Range("K4") "101.6846011131726"
Set-Num $ws.Range("L4") "8.763112482966068"
Set-Num $ws.Range("M4") "0.108572305250173"
Set-Num $ws.Range("N4") "0.1890048347828486"
Set-Num $ws.Range("U4") "0.002621993721129479"
Set-Num $ws.Range("V4") "2.455657407659853"
Set-Num $ws.Range("W4") "0.1176956029219286"

Set-Num $ws.Range("A5") "3"
Set-Num $ws.Range("B5") "4"
Set-Num $ws.Range("D5") "539"
Set-Num $ws.Range("E5") "531"
Set-Num $ws.Range("F5") "7"
Set-Num $ws.Range("G5") "0.4605058979068151"
Set-Num $ws.Range("H5") "64.75077586305163"
Set-Num $ws.Range("I5") "3.55672114852878E-67"
Set-Num $ws.Range("J5") "54.85824257425471"
Set-Num $ws.Range("K5") "101.6846011131726"
Set-Num $ws.Range("L5") "6.689479791273979"
Set-Num $ws.Range("M5") "0.1033111912886153"
Set-Num $ws.Range("N5") "0.1890048347828486"
Set-Num $ws.Range("U5") "0.02960916491904775"
Set-Num $ws.Range("V5") "14.57149069008605"
Set-Num $ws.Range("W5") "6.906177973554183E-07"
